# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table (graphicFrame "Google Shape;122;p17") switches to a
#    different built-in table style.
# 2) The presentation's theme colour scheme is switched from the
#    "Integral / Red Violet" palette to the standard "Office" palette
#    (the font scheme / format scheme were already identical between the
#    two themes, so only the 12 theme colours actually change visually).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{08AA08E7-2510-4E6B-9439-85E20E9C7089}")

# --- 2) Theme colours -------------------------------------------------
# OLE/COM colours are packed as 0x00BBGGRR, i.e. R + G*256 + B*65536.
function ToOleRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches ThemeColorScheme.Colors(1..12):
# Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = ToOleRgb($officeThemeColors[$i - 1])
}
